$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3316.52
$ws.Range("I106").Value = 3153.7058
$ws.Range("J106").Value = 3662.5
$ws.Range("K106").Value = 3153.7058
$ws.Range("L106").Value = 3662.5
$ws.Range("M106").Value = -2522.7058
$ws.Range("N106").Value = -4924.5
$ws.Range("H118").Value = 2469.2856
$ws.Range("I118").Value = 5392.5
$ws.Range("J118").Value = 1300
$ws.Range("K118").Value = 16177.5
$ws.Range("L118").Value = 3900
$ws.Range("M118").Value = -14520.5
$ws.Range("N118").Value = -7214
$ws.Range("H129").Value = 1158592.4
$ws.Range("J129").Value = 1611851.6
$ws.Range("L129").Value = 4835554.800000001
$ws.Range("N129").Value = -4845554.800000001
$ws.Range("H136").Value = 41944.74
$ws.Range("J136").Value = 40552.777
$ws.Range("L136").Value = 40552.777
$ws.Range("N136").Value = -50752.777
$ws.Range("H137").Value = 764.381
$ws.Range("I137").Value = 721.26666
$ws.Range("J137").Value = 788.3333
$ws.Range("K137").Value = 2163.79998
$ws.Range("L137").Value = 2364.9999
$ws.Range("M137").Value = 386.2000200000002
$ws.Range("N137").Value = -7464.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 36000
$ws.Range("J52").Value = 36000
$ws.Range("L52").Value = 36000
$ws.Range("N52").Value = -36636
$ws.Range("H61").Value = 864.42426
$ws.Range("I61").Value = 739.6667
$ws.Range("J61").Value = 1197.1111
$ws.Range("K61").Value = 739.6667
$ws.Range("L61").Value = 1197.1111
$ws.Range("M61").Value = -527.6667
$ws.Range("N61").Value = -1621.1111
$ws.Range("H74").Value = 892.23914
$ws.Range("I74").Value = 855.5227
$ws.Range("J74").Value = 1700
$ws.Range("K74").Value = 855.5227
$ws.Range("L74").Value = 1700
$ws.Range("M74").Value = 18.47730000000001
$ws.Range("N74").Value = -3448
$ws.Range("H77").Value = 892.23914
$ws.Range("I77").Value = 855.5227
$ws.Range("J77").Value = 1700
$ws.Range("K77").Value = 4277.613499999999
$ws.Range("L77").Value = 8500
$ws.Range("M77").Value = 90.38650000000052
$ws.Range("N77").Value = -17236
$ws.Range("H110").Value = 881.7778
$ws.Range("I110").Value = 756
$ws.Range("J110").Value = 1133.3334
$ws.Range("K110").Value = 756
$ws.Range("L110").Value = 1133.3334
$ws.Range("M110").Value = 1289
$ws.Range("N110").Value = -5223.3334
$ws.Range("H128").Value = 35000
$ws.Range("J128").Value = 35000
$ws.Range("L128").Value = 35000
$ws.Range("N128").Value = -44960
$ws.Range("H132").Value = 1315.4193
$ws.Range("I132").Value = 1186
$ws.Range("J132").Value = 1854.6666
$ws.Range("K132").Value = 3558
$ws.Range("L132").Value = 5563.9998
$ws.Range("M132").Value = -1028
$ws.Range("N132").Value = -10623.9998
$ws.Range("H136").Value = 864.42426
$ws.Range("I136").Value = 739.6667
$ws.Range("J136").Value = 1197.1111
$ws.Range("K136").Value = 2219.0001
$ws.Range("L136").Value = 3591.3333
$ws.Range("M136").Value = 330.9998999999998
$ws.Range("N136").Value = -8691.3333

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7193.2607
$ws.Range("I107").Value = 1163.2142
$ws.Range("J107").Value = 16573.334
$ws.Range("K107").Value = 1163.2142
$ws.Range("L107").Value = 16573.334
$ws.Range("M107").Value = 756.7858000000001
$ws.Range("N107").Value = -20413.334
$ws.Range("H139").Value = 52520
$ws.Range("J139").Value = 52520
$ws.Range("L139").Value = 52520
$ws.Range("N139").Value = -62800

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2000
$ws.Range("J4").Value = 2000
$ws.Range("L4").Value = 2000
$ws.Range("N4").Value = -2224
$ws.Range("H16").Value = 1753.8
$ws.Range("I16").Value = 730
$ws.Range("J16").Value = 2777.6
$ws.Range("K16").Value = 730
$ws.Range("L16").Value = 2777.6
$ws.Range("M16").Value = -443
$ws.Range("N16").Value = -3351.6
$ws.Range("H31").Value = 2583.24
$ws.Range("I31").Value = 2379.139
$ws.Range("J31").Value = 3108.0715
$ws.Range("K31").Value = 2379.139
$ws.Range("L31").Value = 3108.0715
$ws.Range("M31").Value = -2084.139
$ws.Range("N31").Value = -3698.0715
$ws.Range("H34").Value = 2583.24
$ws.Range("I34").Value = 2379.139
$ws.Range("J34").Value = 3108.0715
$ws.Range("K34").Value = 2379.139
$ws.Range("L34").Value = 3108.0715
$ws.Range("M34").Value = -2177.139
$ws.Range("N34").Value = -3512.0715
$ws.Range("H113").Value = 1753.8
$ws.Range("I113").Value = 730
$ws.Range("J113").Value = 2777.6
$ws.Range("K113").Value = 730
$ws.Range("L113").Value = 2777.6
$ws.Range("M113").Value = 1440
$ws.Range("N113").Value = -7117.6
$ws.Range("H134").Value = 1196.4849
$ws.Range("I134").Value = 1108.1072
$ws.Range("J134").Value = 1691.4
$ws.Range("K134").Value = 3324.3216
$ws.Range("L134").Value = 5074.200000000001
$ws.Range("M134").Value = -789.3215999999998
$ws.Range("N134").Value = -10144.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 750123.5600000001
$ws.Range("I4").Value = 4040629.8
$ws.Range("J4").Value = 2281.2727
$ws.Range("K4").Value = 12121889.4
$ws.Range("L4").Value = 6843.8181
$ws.Range("M4").Value = -12121777.4
$ws.Range("N4").Value = -7067.8181
$ws.Range("H107").Value = 707786.8
$ws.Range("I107").Value = 801
$ws.Range("J107").Value = 2593082.2
$ws.Range("K107").Value = 2403
$ws.Range("L107").Value = 7779246.600000001
$ws.Range("M107").Value = -483
$ws.Range("N107").Value = -7783086.600000001
$ws.Range("H108").Value = 292.33334
$ws.Range("I108").Value = 292.33334
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 877.0000200000001
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = 2002.99998
$ws.Range("N108").ClearContents()
$ws.Range("H109").Value = 1371.7
$ws.Range("I109").Value = 673.8570999999999
$ws.Range("K109").Value = 2021.5713
$ws.Range("M109").Value = -981.5712999999998
$ws.Range("H110").Value = 5370
$ws.Range("I110").Value = 2962.5
$ws.Range("K110").Value = 8887.5
$ws.Range("M110").Value = -4797.5
$ws.Range("H111").Value = 527
$ws.Range("I111").Value = 527
$ws.Range("K111").Value = 1581
$ws.Range("M111").Value = 1486
$ws.Range("H112").Value = 2712.75
$ws.Range("I112").Value = 1740.6
$ws.Range("J112").Value = 4333
$ws.Range("K112").Value = 5221.799999999999
$ws.Range("L112").Value = 12999
$ws.Range("M112").Value = -4113.799999999999
$ws.Range("N112").Value = -15215
$ws.Range("H115").Value = 1487
$ws.Range("I115").Value = 822.5
$ws.Range("J115").Value = 1930
$ws.Range("K115").Value = 2467.5
$ws.Range("L115").Value = 5790
$ws.Range("M115").Value = -1292.5
$ws.Range("N115").Value = -8140
$ws.Range("H116").Value = 201258
$ws.Range("I116").Value = 1572.5
$ws.Range("K116").Value = 4717.5
$ws.Range("M116").Value = -1275.5
$ws.Range("H117").Value = 3054.1667
$ws.Range("I117").Value = 1000
$ws.Range("J117").Value = 3465
$ws.Range("K117").Value = 3000
$ws.Range("L117").Value = 10395
$ws.Range("M117").Value = 442
$ws.Range("N117").Value = -17279
$ws.Range("H118").Value = 7464.4443
$ws.Range("I118").Value = 743.3333
$ws.Range("J118").Value = 10825
$ws.Range("K118").Value = 2229.9999
$ws.Range("L118").Value = 32475
$ws.Range("M118").Value = -986.9998999999998
$ws.Range("N118").Value = -34961
$ws.Range("H121").Value = 1285.125
$ws.Range("I121").Value = 4000
$ws.Range("J121").Value = 1104.1333
$ws.Range("K121").Value = 12000
$ws.Range("L121").Value = 3312.3999
$ws.Range("M121").Value = -10690
$ws.Range("N121").Value = -5932.3999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4354.077
$ws.Range("I5").Value = 2202.6667
$ws.Range("J5").Value = 4999.5
$ws.Range("K5").Value = 2202.6667
$ws.Range("L5").Value = 4999.5
$ws.Range("M5").Value = -2090.6667
$ws.Range("N5").Value = -5223.5
$ws.Range("H107").Value = 230.27777
$ws.Range("I107").Value = 174.33333
$ws.Range("J107").Value = 342.16666
$ws.Range("K107").Value = 174.33333
$ws.Range("L107").Value = 342.16666
$ws.Range("M107").Value = 1745.66667
$ws.Range("N107").Value = -4182.16666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H46").Value = 2645.2104
$ws.Range("I46").Value = 2046.6666
$ws.Range("J46").Value = 2921.4614
$ws.Range("K46").Value = 2046.6666
$ws.Range("L46").Value = 2921.4614
$ws.Range("M46").Value = -1858.6666
$ws.Range("N46").Value = -3297.4614
$ws.Range("H61").Value = 4465.909
$ws.Range("I61").Value = 5406.8
$ws.Range("J61").Value = 1525.625
$ws.Range("K61").Value = 5406.8
$ws.Range("L61").Value = 1525.625
$ws.Range("M61").Value = -5204.8
$ws.Range("N61").Value = -1929.625
$ws.Range("H113").Value = 4465.909
$ws.Range("I113").Value = 5406.8
$ws.Range("J113").Value = 1525.625
$ws.Range("K113").Value = 5406.8
$ws.Range("L113").Value = 1525.625
$ws.Range("M113").Value = -3236.8
$ws.Range("N113").Value = -5865.625
$ws.Range("H139").Value = 49995
$ws.Range("J139").Value = 49995
$ws.Range("L139").Value = 49995
$ws.Range("N139").Value = -60275

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 159492.31
$ws.Range("J2").Value = 172366.67
$ws.Range("L2").Value = 172366.67
$ws.Range("N2").Value = -172590.67
$ws.Range("H132").Value = 543.76117
$ws.Range("I132").Value = 517.0877
$ws.Range("J132").Value = 695.8
$ws.Range("K132").Value = 1551.2631
$ws.Range("L132").Value = 2087.4
$ws.Range("M132").Value = 978.7368999999999
$ws.Range("N132").Value = -7147.4
$ws.Range("H136").Value = 856.119
$ws.Range("I136").Value = 987
$ws.Range("K136").Value = 2961.0999
$ws.Range("M136").Value = -411
